$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.116.19"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.642.39"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.05"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("E8").Value = "  -2.32%  "
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.59"
$ws.Range("E10").Value = "  -5.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "1.719.86"
$ws.Range("E12").Value = "  +6.31%  "
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.31"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "26.114.34"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "0.0₃0748"
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "189.56"
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.27"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.50"
$ws.Range("E21").Value = "  -4.43%  "
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "144.06"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.76"
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("E27").Value = "  -2.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.21"
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("E30").Value = "  -4.56%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.16"
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.17"
$ws.Range("E32").Value = "  -3.98%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.50"
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("E35").Value = "  -2.93%  "
$ws.Range("D36").Value = "1.122.81"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").Value = "  -5.32%  "
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.73"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.789"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("E42").Value = "  -4.13%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "55.05"
$ws.Range("E44").Value = "  -2.97%  "
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.47"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.416"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.54"
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0924"
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("E51").Value = "  -1.27%  "
